# "Add cantrals by cantons"
# Restructure Sheet1: turn the two-row header (units row removed) into a
# single header row with new column headers (idx, idx2, Name, Date Start,
# Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year),
# and select A2:K2 afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old "units" row (former row 2). This shifts the data rows
#    (former rows 3-14) up to become rows 2-13.
$ws.Rows.Item(2).Delete()

# 2. Build a temporary named cell style carrying just the small Arial-9
#    header font (no explicit number-format flag), apply it to the unit
#    header cells F1:K1, then drop the named style again - this leaves the
#    underlying cell format (a new cellXf entry) in place without adding a
#    permanent named style to the workbook.
$headerStyle = $wb.Styles.Add("TmpHeaderFont")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderFont"
$wb.Styles.Item("TmpHeaderFont").Delete()

# E1 previously carried the old unit-style formatting (style index 1); the
# new layout uses it as a plain header cell with default formatting.
$ws.Range("E1").Style = "Normal"

# 3. Write the new header row text.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# 4. Match the selection recorded in the saved workbook.
$ws.Range("A2:K2").Select() | Out-Null
